$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (year 2025) metrics
$ws.Range("C8").Value = 1053
$ws.Range("E8").Value = 879
$ws.Range("G8").Value = 83.47578347578347
$ws.Range("H8").Value = 16.52421652421652
